$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - randread_128k
$ws.Range("B3").Value = 7062
$ws.Range("C3").Value = 13900
$ws.Range("D3").Value = 22300
$ws.Range("E3").Value = 25800
$ws.Range("F3").Value = 32500
$ws.Range("G3").Value = 29300

# Row 8 - IOPS
$ws.Range("B8").Value = 13300
$ws.Range("C8").Value = 26400
$ws.Range("D8").Value = 52700
$ws.Range("E8").Value = 140000
$ws.Range("F8").Value = 192000
$ws.Range("G8").Value = 181000

# Row 13 - randread_4k
$ws.Range("B13").Value = 4011
$ws.Range("C13").Value = 5067
$ws.Range("D13").Value = 5506
$ws.Range("E13").Value = 5817
$ws.Range("F13").Value = 5673
$ws.Range("G13").Value = 4843

# Row 18 - randwrite_128k
$ws.Range("B18").Value = 117000
$ws.Range("C18").Value = 201000
$ws.Range("D18").Value = 294000
$ws.Range("E18").Value = 175000
$ws.Range("F18").Value = 129000
$ws.Range("G18").Value = 130000

# Row 23 - randwrite_4k
$ws.Range("B23").Value = 3343
$ws.Range("C23").Value = 4895
$ws.Range("D23").Value = 8696
$ws.Range("E23").Value = 11400
$ws.Range("F23").Value = 10000
$ws.Range("G23").Value = 14900

# Row 28 - read_128k
$ws.Range("B28").Value = 112000
$ws.Range("C28").Value = 135000
$ws.Range("D28").Value = 245000
$ws.Range("E28").Value = 344000
$ws.Range("F28").Value = 353000
$ws.Range("G28").Value = 502000

# Row 33 - read_4k
$ws.Range("B33").Value = 2861
$ws.Range("C33").Value = 3155
$ws.Range("D33").Value = 3133
$ws.Range("E33").Value = 3253
$ws.Range("F33").Value = 3392
$ws.Range("G33").Value = 3663

# Row 38 - write_128k
$ws.Range("B38").Value = 99900
$ws.Range("C38").Value = 188000
$ws.Range("D38").Value = 278000
$ws.Range("E38").Value = 401000
$ws.Range("F38").Value = 456000
$ws.Range("G38").Value = 91300
